$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I: pixel size (mm), used to convert DVF errors (supposed isotropic)
# to physically correct values.
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("I1").Font.Bold = $true

$ws.Range("I2").Value = 1.818

$ws.Range("I5").Select()
